# Append the 2025-03-26 price row to every "Solar_Prices" sheet.
# Each sheet has a Date/Price table in columns A:B ending at row 24
# (the most recent date, 2025-03-25). We add a new row 25 for
# 2025-03-26, carrying forward the same price as the prior day
# (row 24), and let the sheet's dimension/used range grow to B25.

$wb = $excel.ActiveWorkbook
$newDate = "2025-03-26"

for ($i = 1; $i -le $wb.Worksheets.Count(); $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Find the last used row in column A/B (currently row 24 for all sheets).
    $lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()
    $newRow = $lastRow + 1

    $priorPrice = $ws.Cells.Item($lastRow, 2).Value()

    # Force the new cells to be stored as text (matching the existing
    # inline/shared-string "Date"/"Price" columns) rather than letting
    # Excel auto-convert the date-like string into a date serial number
    # or the numeric-looking price into a real number.
    $ws.Cells.Item($newRow, 1).Value = "'" + $newDate
    $ws.Cells.Item($newRow, 2).Value = "'" + $priorPrice
}
